# Updates cryptos list values (Price / Volume(1h)) per the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.503.54'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').Value = '3.051.38'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'517.55"
$ws.Range('E5').Value = '  +2.54%  '
$ws.Range('D6').Value = "'141.18"
$ws.Range('E6').Value = '  +2.95%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('D9').Value = "'7.26"
$ws.Range('E9').Value = '  +1.14%  '
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('D11').Value = "'0.375"
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('D12').Value = '3.578.17'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('E13').Value = '  +3.27%  '
$ws.Range('D14').Value = "'25.51"
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('D16').Value = '57.544.74'
$ws.Range('D17').Value = '3.051.88'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = "'12.78"
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').Value = "'8.10"
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = "'329.53"
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').Value = "'65.71"
$ws.Range('E24').Value = '  +1.88%  '
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '0.0₃0895'
$ws.Range('E27').Value = '  -2.56%  '
$ws.Range('D28').Value = "'6.30"
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').Value = "'7.13"
$ws.Range('E29').Value = '  +1.98%  '
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('D31').Value = "'20.61"
$ws.Range('E32').Value = '  +1.79%  '
$ws.Range('D33').Value = "'154.50"
$ws.Range('E33').Value = '  +1.18%  '
$ws.Range('D34').Value = "'27.26"
$ws.Range('E34').Value = '  +5.74%  '
$ws.Range('D35').Value = "'4.47"
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('E36').Value = '  +2.27%  '
$ws.Range('E37').Value = '  +1.61%  '
$ws.Range('D38').Value = "'0.0671"
$ws.Range('E38').Value = '  +1.70%  '
$ws.Range('D39').Value = '3.089.81'
$ws.Range('E39').Value = '  +2.47%  '
$ws.Range('D40').Value = "'3.89"
$ws.Range('E40').Value = '  +2.91%  '
$ws.Range('D41').Value = "'36.63"
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').Value = "'0.649"
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '2.248.47'
$ws.Range('E44').Value = '  +3.61%  '
$ws.Range('D45').Value = "'0.0257"
$ws.Range('E45').Value = '  +8.78%  '
$ws.Range('D46').Value = "'20.68"
$ws.Range('E46').Value = '  +6.45%  '
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('D48').Value = "'5.85"
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').Value = "'0.915"
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').Value = "'260.20"
$ws.Range('E50').Value = '  +15.01%  '
$ws.Range('D51').Value = "'0.712"
$ws.Range('E51').Value = '  +6.32%  '
